$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Insert 20 new rows before the existing data (old rows 2-21 shift down to 22-41)
$ws.Range("A2:A21").EntireRow.Insert()

# Populate the newly inserted rows 2-21 with the new player performance data
$ws.Cells.Item(2,1).Value = "'3458"
$ws.Cells.Item(2,2).Value = 3
$ws.Cells.Item(2,3).Value = "'0"
$ws.Cells.Item(2,4).Value = "'0"
$ws.Cells.Item(2,5).Value = "'0.98%"
$ws.Cells.Item(2,6).Value = "NO"
$ws.Cells.Item(3,1).Value = "'3471"
$ws.Cells.Item(3,2).Value = 3
$ws.Cells.Item(3,3).Value = "'1"
$ws.Cells.Item(3,4).Value = "'0"
$ws.Cells.Item(3,5).Value = "'11.27%"
$ws.Cells.Item(3,6).Value = "NO"
$ws.Cells.Item(4,1).Value = "'3472"
$ws.Cells.Item(4,2).Value = 2
$ws.Cells.Item(4,3).Value = "'0"
$ws.Cells.Item(4,4).Value = "'0"
$ws.Cells.Item(4,5).Value = "'1.13%"
$ws.Cells.Item(4,6).Value = "NO"
$ws.Cells.Item(5,1).Value = "'3886"
$ws.Cells.Item(5,2).Value = 1
$ws.Cells.Item(5,3).Value = "'7"
$ws.Cells.Item(5,4).Value = "'0"
$ws.Cells.Item(5,5).Value = "'17.67%"
$ws.Cells.Item(5,6).Value = "NO"
$ws.Cells.Item(6,1).Value = "'3888"
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = "'4"
$ws.Cells.Item(6,4).Value = "'2"
$ws.Cells.Item(6,5).Value = "'23.04%"
$ws.Cells.Item(6,6).Value = "NO"
$ws.Cells.Item(7,1).Value = "'3893"
$ws.Cells.Item(7,2).Value = "'"
$ws.Cells.Item(7,3).Value = "'"
$ws.Cells.Item(7,4).Value = "'"
$ws.Cells.Item(7,5).Value = "'"
$ws.Cells.Item(7,6).Value = "NO"
$ws.Cells.Item(8,1).Value = "'3894"
$ws.Cells.Item(8,2).Value = 3
$ws.Cells.Item(8,3).Value = "'0"
$ws.Cells.Item(8,4).Value = "'0"
$ws.Cells.Item(8,5).Value = "'1.41%"
$ws.Cells.Item(8,6).Value = "NO"
$ws.Cells.Item(9,1).Value = "'3896"
$ws.Cells.Item(9,2).Value = 3
$ws.Cells.Item(9,3).Value = "'4"
$ws.Cells.Item(9,4).Value = "'1"
$ws.Cells.Item(9,5).Value = "'20.49%"
$ws.Cells.Item(9,6).Value = "NO"
$ws.Cells.Item(10,1).Value = "'3898"
$ws.Cells.Item(10,2).Value = 1
$ws.Cells.Item(10,3).Value = "'4"
$ws.Cells.Item(10,4).Value = "'3"
$ws.Cells.Item(10,5).Value = "'36.98%"
$ws.Cells.Item(10,6).Value = "NO"
$ws.Cells.Item(11,1).Value = "'3903"
$ws.Cells.Item(11,2).Value = "'"
$ws.Cells.Item(11,3).Value = "'"
$ws.Cells.Item(11,4).Value = "'"
$ws.Cells.Item(11,5).Value = "'"
$ws.Cells.Item(11,6).Value = "NO"
$ws.Cells.Item(12,1).Value = "'3905"
$ws.Cells.Item(12,2).Value = "'"
$ws.Cells.Item(12,3).Value = "'"
$ws.Cells.Item(12,4).Value = "'"
$ws.Cells.Item(12,5).Value = "'"
$ws.Cells.Item(12,6).Value = "NO"
$ws.Cells.Item(13,1).Value = "'3909"
$ws.Cells.Item(13,2).Value = 1
$ws.Cells.Item(13,3).Value = "'3"
$ws.Cells.Item(13,4).Value = "'0"
$ws.Cells.Item(13,5).Value = "'5.19%"
$ws.Cells.Item(13,6).Value = "NO"
$ws.Cells.Item(14,1).Value = "'3929"
$ws.Cells.Item(14,2).Value = "'"
$ws.Cells.Item(14,3).Value = "'"
$ws.Cells.Item(14,4).Value = "'"
$ws.Cells.Item(14,5).Value = "'"
$ws.Cells.Item(14,6).Value = "NO"
$ws.Cells.Item(15,1).Value = "'3931"
$ws.Cells.Item(15,2).Value = 3
$ws.Cells.Item(15,3).Value = "'1"
$ws.Cells.Item(15,4).Value = "'0"
$ws.Cells.Item(15,5).Value = "'3.02%"
$ws.Cells.Item(15,6).Value = "NO"
$ws.Cells.Item(16,1).Value = "'3937"
$ws.Cells.Item(16,2).Value = 1
$ws.Cells.Item(16,3).Value = "'8"
$ws.Cells.Item(16,4).Value = "'2"
$ws.Cells.Item(16,5).Value = "'41.21%"
$ws.Cells.Item(16,6).Value = "YES"
$ws.Cells.Item(17,1).Value = "'3973"
$ws.Cells.Item(17,2).Value = "'"
$ws.Cells.Item(17,3).Value = "'"
$ws.Cells.Item(17,4).Value = "'"
$ws.Cells.Item(17,5).Value = "'"
$ws.Cells.Item(17,6).Value = "NO"
$ws.Cells.Item(18,1).Value = "'3975"
$ws.Cells.Item(18,2).Value = 2
$ws.Cells.Item(18,3).Value = "'2"
$ws.Cells.Item(18,4).Value = "'0"
$ws.Cells.Item(18,5).Value = "'3.40%"
$ws.Cells.Item(18,6).Value = "NO"
$ws.Cells.Item(19,1).Value = "'3977"
$ws.Cells.Item(19,2).Value = 1
$ws.Cells.Item(19,3).Value = "'2"
$ws.Cells.Item(19,4).Value = "'0"
$ws.Cells.Item(19,5).Value = "'8.50%"
$ws.Cells.Item(19,6).Value = "NO"
$ws.Cells.Item(20,1).Value = "'4234"
$ws.Cells.Item(20,2).Value = 3
$ws.Cells.Item(20,3).Value = "'6"
$ws.Cells.Item(20,4).Value = "'0"
$ws.Cells.Item(20,5).Value = "'20.49%"
$ws.Cells.Item(20,6).Value = "NO"
$ws.Cells.Item(21,1).Value = "'4235"
$ws.Cells.Item(21,2).Value = "'"
$ws.Cells.Item(21,3).Value = "'"
$ws.Cells.Item(21,4).Value = "'"
$ws.Cells.Item(21,5).Value = "'"
$ws.Cells.Item(21,6).Value = "NO"
